# Generate Report for Handback
#
# The handback file "6e192eb7-ea71-4498-b5b3-74c6c8e2fb1e.md" has now been
# handed back and is in sync with en-US, so refresh the localization status
# report: flip its status from "Ready for handoff" to
# "Handed back: in sync with en-US", stamp the new handback datetimes, and
# clear the now-stale "version mismatch" error detail for each language.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for the 6e192eb7 row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# --- zh-cn sheet: Status / Latest Handback DateTime / Error Detail ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("K3").Value = "2016-08-22 06:47:33"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8

# --- de-de sheet: Status / Latest Handback DateTime / Error Detail ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("K3").Value = "2016-08-22 06:47:40"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8
